$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 135
$ws.Range("D21").Value = 124
$ws.Range("E21").Value = 11
$ws.Range("F21").Value = 35.53008595988539
